$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3105.8333
$ws.Range("I62").Value = 3105.8333
$ws.Range("K62").Value = 3105.8333
$ws.Range("M62").Value = -2481.8333

$ws.Range("H64").Value = 3955
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 3932.5
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 3932.5
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -4428.5

$ws.Range("H65").Value = 3105.8333
$ws.Range("I65").Value = 3105.8333
$ws.Range("K65").Value = 15529.1665
$ws.Range("M65").Value = -12409.1665

$ws.Range("H67").Value = 3955
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 3932.5
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 3932.5
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -5648.5

$ws.Range("H70").Value = 1109806.5

$ws.Range("H73").Value = 1109806.5

$ws.Range("H99").Value = 3478.3333
$ws.Range("I99").Value = 159
$ws.Range("J99").Value = 3893.25
$ws.Range("K99").Value = 477
$ws.Range("L99").Value = 11679.75
$ws.Range("M99").Value = 1021
$ws.Range("N99").Value = -14675.75

$ws.Range("H111").Value = 6999.25
$ws.Range("I111").Value = 1998.5
$ws.Range("K111").Value = 5995.5
$ws.Range("M111").Value = -2928.5

$ws.Range("H132").Value = 4504.4185
$ws.Range("I132").Value = 2964
$ws.Range("K132").Value = 8892
$ws.Range("M132").Value = -6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2143.3823
$ws.Range("I2").Value = 2243.2693
$ws.Range("K2").Value = 2243.2693
$ws.Range("M2").Value = -2130.2693

$ws.Range("H5").Value = 297.5
$ws.Range("I5").Value = 40
$ws.Range("K5").Value = 40
$ws.Range("M5").Value = 72

$ws.Range("H74").Value = 4415.5454
$ws.Range("I74").Value = 4034.625
$ws.Range("J74").Value = 5431.3335
$ws.Range("K74").Value = 4034.625
$ws.Range("L74").Value = 5431.3335
$ws.Range("M74").Value = -3160.625
$ws.Range("N74").Value = -7179.3335

$ws.Range("H77").Value = 4415.5454
$ws.Range("I77").Value = 4034.625
$ws.Range("J77").Value = 5431.3335
$ws.Range("K77").Value = 20173.125
$ws.Range("L77").Value = 27156.6675
$ws.Range("M77").Value = -15805.125
$ws.Range("N77").Value = -35892.6675

$ws.Range("H110").Value = 6912
$ws.Range("I110").Value = 8802.5
$ws.Range("J110").Value = 3131
$ws.Range("K110").Value = 8802.5
$ws.Range("L110").Value = 3131
$ws.Range("M110").Value = -6757.5
$ws.Range("N110").Value = -7221

$ws.Range("H116").Value = 2143.3823
$ws.Range("I116").Value = 2243.2693
$ws.Range("K116").Value = 2243.2693
$ws.Range("M116").Value = 50.73070000000007

$ws.Range("H122").Value = 3692.7
$ws.Range("I122").Value = 3887.611
$ws.Range("J122").Value = 1938.5
$ws.Range("K122").Value = 11662.833
$ws.Range("L122").Value = 5815.5
$ws.Range("M122").Value = -9212.832999999999
$ws.Range("N122").Value = -10715.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2143.3823
$ws.Range("I3").Value = 2243.2693
$ws.Range("K3").Value = 2243.2693
$ws.Range("M3").Value = -2129.2693

$ws.Range("H4").Value = 297.5
$ws.Range("I4").Value = 40
$ws.Range("K4").Value = 40
$ws.Range("M4").Value = 75

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 1139.6666
$ws.Range("I22").Value = 149.75
$ws.Range("K22").Value = 149.75
$ws.Range("M22").Value = 23.25

$ws.Range("H99").Value = 2190.818
$ws.Range("I99").Value = 2259.9
$ws.Range("K99").Value = 2259.9
$ws.Range("M99").Value = -761.9000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 450
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H31").Value = 45457784
$ws.Range("J31").Value = 2937
$ws.Range("L31").Value = 2937
$ws.Range("N31").Value = -3527

$ws.Range("H34").Value = 45457784
$ws.Range("J34").Value = 2937
$ws.Range("L34").Value = 2937
$ws.Range("N34").Value = -3341

$ws.Range("H87").Value = 99999
$ws.Range("J87").Value = 99999
$ws.Range("L87").Value = 99999
$ws.Range("N87").Value = -102371

$ws.Range("H90").Value = 99999
$ws.Range("J90").Value = 99999
$ws.Range("L90").Value = 299997
$ws.Range("N90").Value = -311853

$ws.Range("H122").Value = 4335.5884
$ws.Range("I122").Value = 4396.222
$ws.Range("J122").Value = 4267.375
$ws.Range("K122").Value = 13188.666
$ws.Range("L122").Value = 12802.125
$ws.Range("M122").Value = -10738.666
$ws.Range("N122").Value = -17702.125

$ws.Range("H132").Value = 2941.04
$ws.Range("I132").Value = 2682
$ws.Range("J132").Value = 3401.5557
$ws.Range("K132").Value = 8046
$ws.Range("L132").Value = 10204.6671
$ws.Range("M132").Value = -5516
$ws.Range("N132").Value = -15264.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 121
$ws.Range("I40").Value = 45.2
$ws.Range("K40").Value = 180.8
$ws.Range("M40").Value = -111.8

$ws.Range("H122").Value = 57141.168
$ws.Range("J122").Value = 2416.3333
$ws.Range("L122").Value = 21746.9997
$ws.Range("N122").Value = -26646.9997

$ws.Range("H137").Value = 3706.7273
$ws.Range("I137").Value = 1436.9
$ws.Range("J137").Value = 5598.25
$ws.Range("K137").Value = 4310.700000000001
$ws.Range("L137").Value = 16794.75
$ws.Range("M137").Value = 789.2999999999993
$ws.Range("N137").Value = -26994.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1545039.9
$ws.Range("I113").Value = 1784.8889
$ws.Range("K113").Value = 1784.8889
$ws.Range("M113").Value = 385.1111000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2328.4285
$ws.Range("I46").Value = 699.5
$ws.Range("J46").Value = 2980
$ws.Range("K46").Value = 699.5
$ws.Range("L46").Value = 2980
$ws.Range("M46").Value = -511.5
$ws.Range("N46").Value = -3356

$ws.Range("H61").Value = 4999
$ws.Range("J61").Value = 4999
$ws.Range("L61").Value = 4999
$ws.Range("N61").Value = -5403

$ws.Range("H68").Value = 3790137
$ws.Range("I68").Value = 8334940
$ws.Range("J68").Value = 2800.8333
$ws.Range("K68").Value = 8334940
$ws.Range("L68").Value = 2800.8333
$ws.Range("M68").Value = -8334191
$ws.Range("N68").Value = -4298.8333

$ws.Range("H71").Value = 3790137
$ws.Range("I71").Value = 8334940
$ws.Range("J71").Value = 2800.8333
$ws.Range("K71").Value = 41674700
$ws.Range("L71").Value = 14004.1665
$ws.Range("M71").Value = -41670956
$ws.Range("N71").Value = -21492.1665

$ws.Range("H93").Value = 1854960.4
$ws.Range("J93").Value = 5056722
$ws.Range("L93").Value = 5056722
$ws.Range("N93").Value = -5059218

$ws.Range("H113").Value = 4999
$ws.Range("J113").Value = 4999
$ws.Range("L113").Value = 4999
$ws.Range("N113").Value = -9339

$ws.Range("H132").Value = 2237.28
$ws.Range("I132").Value = 2082.8125
$ws.Range("J132").Value = 2511.889
$ws.Range("K132").Value = 6248.4375
$ws.Range("L132").Value = 7535.667
$ws.Range("M132").Value = -3718.4375
$ws.Range("N132").Value = -12595.667

$ws.Range("H136").Value = 2908.24
$ws.Range("J136").Value = 2934.1667
$ws.Range("L136").Value = 8802.500100000001
$ws.Range("N136").Value = -13902.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 10261.111
$ws.Range("I96").Value = 9228.333000000001
$ws.Range("J96").Value = 12326.667
$ws.Range("K96").Value = 9228.333000000001
$ws.Range("L96").Value = 12326.667
$ws.Range("M96").Value = -7855.333000000001
$ws.Range("N96").Value = -15072.667

$ws.Range("H119").Value = 99999
$ws.Range("J119").Value = 99999
$ws.Range("L119").Value = 99999
$ws.Range("N119").Value = -109675
